$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.892.74"
$ws.Range("E2").Value = "  -4.75%  "
$ws.Range("D3").Value = "2.600.76"
$ws.Range("E3").Value = "  -3.58%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'514.87"
$ws.Range("E5").Value = "  -2.15%  "
$ws.Range("D6").Value = "'141.42"
$ws.Range("E6").Value = "  -2.54%  "
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("D8").Value = "'0.565"
$ws.Range("E8").Value = "  -1.82%  "
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("E10").Value = "  -3.14%  "
$ws.Range("E11").Value = "  -1.02%  "
$ws.Range("E12").Value = "  +1.10%  "
$ws.Range("D13").Value = "3.065.36"
$ws.Range("E13").Value = "  -3.61%  "
$ws.Range("D14").Value = "57.899.18"
$ws.Range("E14").Value = "  -4.64%  "
$ws.Range("D15").Value = "'20.57"
$ws.Range("E15").Value = "  -3.22%  "
$ws.Range("E16").Value = "  -2.00%  "
$ws.Range("D17").Value = "2.617.39"
$ws.Range("E17").Value = "  -3.84%  "
$ws.Range("D18").Value = "'4.38"
$ws.Range("E18").Value = "  -2.67%  "
$ws.Range("D19").Value = "'333.17"
$ws.Range("E19").Value = "  -3.51%  "
$ws.Range("D20").Value = "'10.28"
$ws.Range("E20").Value = "  -2.98%  "
$ws.Range("D21").Value = "'6.22"
$ws.Range("E21").Value = "  -3.26%  "
$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("D23").Value = "'63.82"
$ws.Range("E23").Value = "  +0.35%  "
$ws.Range("D24").Value = "'0.414"
$ws.Range("E24").Value = "  -1.61%  "
$ws.Range("D25").Value = "'0.165"
$ws.Range("E25").Value = "  -2.71%  "
$ws.Range("E26").Value = "  +0.65%  "
$ws.Range("E27").Value = "  -3.46%  "
$ws.Range("E28").Value = "  -4.49%  "
$ws.Range("D29").Value = "'6.57"
$ws.Range("E29").Value = "  -3.90%  "
$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "  +0.13%  "
$ws.Range("E31").Value = "  -1.81%  "
$ws.Range("D32").Value = "'150.83"
$ws.Range("E32").Value = "  +0.46%  "
$ws.Range("D33").Value = "'18.60"
$ws.Range("E33").Value = "  -2.31%  "
$ws.Range("E34").Value = "  -4.11%  "
$ws.Range("E35").Value = "  -5.94%  "
$ws.Range("D36").Value = "'0.890"
$ws.Range("E36").Value = "  -5.44%  "
$ws.Range("D37").Value = "'36.46"
$ws.Range("E37").Value = "  -1.85%  "
$ws.Range("D38").Value = "'0.836"
$ws.Range("E38").Value = "  -4.11%  "
$ws.Range("D39").Value = "'1.42"
$ws.Range("E39").Value = "  -6.59%  "
$ws.Range("D40").Value = "'3.59"
$ws.Range("E40").Value = "  -2.11%  "
$ws.Range("D41").Value = "'0.999"
$ws.Range("E41").Value = "  +0.22%  "
$ws.Range("D42").Value = "'0.0961"
$ws.Range("E42").Value = "  -2.45%  "
$ws.Range("D44").Value = "'266.88"
$ws.Range("E44").Value = "  -5.47%  "
$ws.Range("E45").Value = "  +1.14%  "
$ws.Range("D46").Value = "'18.98"
$ws.Range("E46").Value = "  -5.45%  "
$ws.Range("E47").Value = "  -1.71%  "
$ws.Range("D48").Value = "2.029.88"
$ws.Range("E48").Value = "  -5.33%  "
$ws.Range("D49").Value = "'0.0227"
$ws.Range("E49").Value = "  -2.33%  "
$ws.Range("D50").Value = "'4.58"
$ws.Range("E50").Value = "  -4.79%  "
$ws.Range("D51").Value = "'18.16"
$ws.Range("E51").Value = "  -4.92%  "
